$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- "About" sheet updates ---
$aboutWs = $wb.Worksheets.Item("About")

$aboutWs.Range("A2").Value = "Version: $newVersion"

$aboutWs.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Junde Coal Mine, China, M1156, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet updates ---
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 12; $row++) {
    $cell = $dataWs.Cells.Item($row, 19)  # column S
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
